$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 10890354
$ws.Range("I12").Value = 17424492
$ws.Range("K12").Value = 17424492
$ws.Range("M12").Value = -17424322
$ws.Range("H48").Value = 3000
$ws.Range("J48").Value = 3000
$ws.Range("L48").Value = 9000
$ws.Range("N48").Value = -9584
$ws.Range("H56").Value = 3000
$ws.Range("J56").Value = 3000
$ws.Range("L56").Value = 9000
$ws.Range("N56").Value = -10068
$ws.Range("H62").Value = 2285
$ws.Range("I62").Value = 2285
$ws.Range("K62").Value = 2285
$ws.Range("M62").Value = -1661
$ws.Range("H65").Value = 2285
$ws.Range("I65").Value = 2285
$ws.Range("K65").Value = 11425
$ws.Range("M65").Value = -8305
$ws.Range("H112").Value = 4110.32
$ws.Range("J112").Value = 4421.8096
$ws.Range("L112").Value = 13265.4288
$ws.Range("N112").Value = -15481.4288
$ws.Range("H132").Value = 26236.95
$ws.Range("I132").Value = 27459.947
$ws.Range("K132").Value = 82379.841
$ws.Range("M132").Value = -79849.841
$ws.Range("H137").Value = 12501324
$ws.Range("I137").Value = 1224.5
$ws.Range("J137").Value = 25001422
$ws.Range("K137").Value = 3673.5
$ws.Range("L137").Value = 75004266
$ws.Range("M137").Value = -1123.5
$ws.Range("N137").Value = -75009366
$ws.Range("H138").Value = 5043.8906
$ws.Range("I138").Value = 9979.571
$ws.Range("J138").Value = 3661.9
$ws.Range("K138").Value = 29938.713
$ws.Range("L138").Value = 10985.7
$ws.Range("M138").Value = -24798.713
$ws.Range("N138").Value = -21265.7

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1108.7916
$ws.Range("I2").Value = 734.619
$ws.Range("J2").Value = 3728
$ws.Range("K2").Value = 734.619
$ws.Range("L2").Value = 3728
$ws.Range("M2").Value = -621.619
$ws.Range("N2").Value = -3954
$ws.Range("H32").Value = 148293.2
$ws.Range("I32").Value = 159450.97
$ws.Range("J32").Value = 7705.4
$ws.Range("K32").Value = 159450.97
$ws.Range("L32").Value = 7705.4
$ws.Range("M32").Value = -159163.97
$ws.Range("N32").Value = -8279.4
$ws.Range("H38").Value = 0
$ws.Range("I38").Value = 0
$ws.Range("K38").Value = 0
$ws.Range("M38").ClearContents()
$ws.Range("H68").Value = 25000
$ws.Range("J68").Value = 25000
$ws.Range("L68").Value = 25000
$ws.Range("N68").Value = -26622
$ws.Range("H71").Value = 25000
$ws.Range("J71").Value = 25000
$ws.Range("L71").Value = 75000
$ws.Range("N71").Value = -83112
$ws.Range("H97").Value = 8140.9287
$ws.Range("I97").Value = 10115.818
$ws.Range("K97").Value = 10115.818
$ws.Range("M97").Value = -9619.817999999999
$ws.Range("H116").Value = 1108.7916
$ws.Range("I116").Value = 734.619
$ws.Range("J116").Value = 3728
$ws.Range("K116").Value = 734.619
$ws.Range("L116").Value = 3728
$ws.Range("M116").Value = 1559.381
$ws.Range("N116").Value = -8316

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1108.7916
$ws.Range("I3").Value = 734.619
$ws.Range("J3").Value = 3728
$ws.Range("K3").Value = 734.619
$ws.Range("L3").Value = 3728
$ws.Range("M3").Value = -620.619
$ws.Range("N3").Value = -3956
$ws.Range("H20").Value = 968.73334
$ws.Range("I20").Value = 855.4761999999999
$ws.Range("K20").Value = 855.4761999999999
$ws.Range("M20").Value = -608.4761999999999
$ws.Range("H94").Value = 5926.067
$ws.Range("I94").Value = 4989.4546
$ws.Range("K94").Value = 4989.4546
$ws.Range("M94").Value = -4538.4546
$ws.Range("H134").Value = 18751816
$ws.Range("I134").Value = 1692.3243
$ws.Range("J134").Value = 81820420
$ws.Range("K134").Value = 5076.9729
$ws.Range("L134").Value = 245461260
$ws.Range("M134").Value = -2541.9729
$ws.Range("N134").Value = -245466330

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3547
$ws.Range("I16").Value = 3127.2307
$ws.Range("J16").Value = 4911.25
$ws.Range("K16").Value = 3127.2307
$ws.Range("L16").Value = 4911.25
$ws.Range("M16").Value = -2840.2307
$ws.Range("N16").Value = -5485.25
$ws.Range("H113").Value = 3547
$ws.Range("I113").Value = 3127.2307
$ws.Range("J113").Value = 4911.25
$ws.Range("K113").Value = 3127.2307
$ws.Range("L113").Value = 4911.25
$ws.Range("M113").Value = -957.2307000000001
$ws.Range("N113").Value = -9251.25
$ws.Range("H132").Value = 24730.4
$ws.Range("I132").Value = 28641.71
$ws.Range("J132").Value = 3497.5715
$ws.Range("K132").Value = 85925.13
$ws.Range("L132").Value = 10492.7145
$ws.Range("M132").Value = -83395.13
$ws.Range("N132").Value = -15552.7145

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 309.09525
$ws.Range("I2").Value = 215.42857
$ws.Range("J2").Value = 496.42856
$ws.Range("K2").Value = 1292.57142
$ws.Range("L2").Value = 2978.57136
$ws.Range("M2").Value = -1179.57142
$ws.Range("N2").Value = -3204.57136
$ws.Range("H8").Value = 1877
$ws.Range("I8").Value = 1877
$ws.Range("K8").Value = 5631
$ws.Range("M8").Value = -5492
$ws.Range("H69").Value = 9697.177
$ws.Range("I69").Value = 27995
$ws.Range("K69").Value = 83985
$ws.Range("M69").Value = -83174
$ws.Range("H72").Value = 9697.177
$ws.Range("I72").Value = 27995
$ws.Range("K72").Value = 251955
$ws.Range("M72").Value = -247899
$ws.Range("H75").Value = 5332.5
$ws.Range("J75").Value = 6999.25
$ws.Range("L75").Value = 20997.75
$ws.Range("N75").Value = -22993.75
$ws.Range("H78").Value = 5332.5
$ws.Range("J78").Value = 6999.25
$ws.Range("L78").Value = 62993.25
$ws.Range("N78").Value = -72977.25
$ws.Range("H86").Value = 474.45456
$ws.Range("I86").Value = 68
$ws.Range("J86").Value = 564.7778
$ws.Range("K86").Value = 204
$ws.Range("L86").Value = 1694.3334
$ws.Range("M86").Value = 982
$ws.Range("N86").Value = -4066.3334
$ws.Range("H89").Value = 474.45456
$ws.Range("I89").Value = 68
$ws.Range("J89").Value = 564.7778
$ws.Range("K89").Value = 612
$ws.Range("L89").Value = 5083.000199999999
$ws.Range("M89").Value = 5316
$ws.Range("N89").Value = -16939.0002
$ws.Range("H121").Value = 1093.2858
$ws.Range("J121").Value = 2853
$ws.Range("L121").Value = 8559
$ws.Range("N121").Value = -11179
$ws.Range("H122").Value = 5322903.5
$ws.Range("J122").Value = 2381969
$ws.Range("L122").Value = 21437721
$ws.Range("N122").Value = -21442621
$ws.Range("H140").Value = 1961.5
$ws.Range("I140").Value = 1594.409
$ws.Range("J140").Value = 5999.5
$ws.Range("K140").Value = 4783.227000000001
$ws.Range("L140").Value = 17998.5
$ws.Range("M140").Value = 396.7729999999992
$ws.Range("N140").Value = -28358.5

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H49").Value = 45999
$ws.Range("J49").Value = 45999
$ws.Range("L49").Value = 45999
$ws.Range("N49").Value = -46367
$ws.Range("H70").Value = 3999.5
$ws.Range("I70").Value = 3999.5
$ws.Range("K70").Value = 3999.5
$ws.Range("M70").Value = -3729.5
$ws.Range("H73").Value = 3999.5
$ws.Range("I73").Value = 3999.5
$ws.Range("K73").Value = 3999.5
$ws.Range("M73").Value = -3063.5
$ws.Range("H113").Value = 1088.3
$ws.Range("J113").Value = 1147
$ws.Range("L113").Value = 1147
$ws.Range("N113").Value = -5487
$ws.Range("H122").Value = 2731.3667
$ws.Range("J122").Value = 3819
$ws.Range("L122").Value = 11457
$ws.Range("N122").Value = -16357

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()
$ws.Range("H22").Value = 4606.628
$ws.Range("I22").Value = 2454.5
$ws.Range("K22").Value = 2454.5
$ws.Range("M22").Value = -2159.5
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("M25").ClearContents()
$ws.Range("H27").Value = 4606.628
$ws.Range("I27").Value = 2454.5
$ws.Range("K27").Value = 2454.5
$ws.Range("M27").Value = -2347.5
$ws.Range("H46").Value = 12999.786
$ws.Range("J46").Value = 7749.8335
$ws.Range("L46").Value = 7749.8335
$ws.Range("N46").Value = -8125.8335
$ws.Range("H56").Value = 9999
$ws.Range("I56").Value = 9999
$ws.Range("K56").Value = 9999
$ws.Range("M56").Value = -9308
$ws.Range("H61").Value = 4203.0386
$ws.Range("I61").Value = 4914.0713
$ws.Range("J61").Value = 3373.5
$ws.Range("K61").Value = 4914.0713
$ws.Range("L61").Value = 3373.5
$ws.Range("M61").Value = -4712.0713
$ws.Range("N61").Value = -3777.5
$ws.Range("H93").Value = 3214.3333
$ws.Range("I93").Value = 2404.6667
$ws.Range("J93").Value = 4024
$ws.Range("K93").Value = 2404.6667
$ws.Range("L93").Value = 4024
$ws.Range("M93").Value = -1156.6667
$ws.Range("N93").Value = -6520
$ws.Range("H100").Value = 2487
$ws.Range("I100").Value = 2217.9375
$ws.Range("J100").Value = 2917.5
$ws.Range("K100").Value = 2217.9375
$ws.Range("L100").Value = 2917.5
$ws.Range("M100").Value = -1676.9375
$ws.Range("N100").Value = -3999.5
$ws.Range("H113").Value = 4203.0386
$ws.Range("I113").Value = 4914.0713
$ws.Range("J113").Value = 3373.5
$ws.Range("K113").Value = 4914.0713
$ws.Range("L113").Value = 3373.5
$ws.Range("M113").Value = -2744.0713
$ws.Range("N113").Value = -7713.5
$ws.Range("H132").Value = 2789.318
$ws.Range("I132").Value = 2775.3823
$ws.Range("J132").Value = 2836.7
$ws.Range("K132").Value = 8326.1469
$ws.Range("L132").Value = 8510.099999999999
$ws.Range("M132").Value = -5796.1469
$ws.Range("N132").Value = -13570.1

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H95").Value = 36249.75
$ws.Range("J95").Value = 39999.668
$ws.Range("L95").Value = 39999.668
$ws.Range("N95").Value = -45491.668
$ws.Range("H122").Value = 1363.7709
$ws.Range("I122").Value = 1163.0465
$ws.Range("J122").Value = 3090
$ws.Range("K122").Value = 3489.1395
$ws.Range("L122").Value = 9270
$ws.Range("M122").Value = -1039.1395
$ws.Range("N122").Value = -14170
$ws.Range("H132").Value = 2518.48
$ws.Range("I132").Value = 1829.2106
$ws.Range("K132").Value = 5487.6318
$ws.Range("M132").Value = -2957.6318
